$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.950247585773468
$ws.Range("B1").Value = 2.075462102890015
$ws.Range("C1").Value = 8.086058616638184
$ws.Range("D1").Value = 2.322593212127686
$ws.Range("E1").Value = 0.9029766917228699
